$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.298.99"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.082.25"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'327.56"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5202"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D8").Value = "'0.4306"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "'0.08830"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("D11").Value = "'1.160"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("D12").Value = "'24.46"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "2.083.45"
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").Value = "'6.712"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").Value = "'7.672"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "'95.20"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'0.00001120"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'0.06625"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "'18.84"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'6.304"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "30.337.00"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  +3.92%  "
$ws.Range("D25").Value = "'2.290"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").Value = "2.328.22"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "'22.30"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'2.580"
$ws.Range("E28").Value = "  +5.78%  "
$ws.Range("D29").Value = "'162.04"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "'131.01"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "'1.190"
$ws.Range("E31").Value = "  +3.93%  "
$ws.Range("D32").Value = "'0.1067"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  +20.09%  "
$ws.Range("D34").Value = "'6.180"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "'3.822"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02581"
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'9.836"
$ws.Range("E37").Value = "  +7.95%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06681"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'12.67"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").Value = "'5.443"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "'0.2252"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").Value = "'0.6811"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'14.03"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "'0.6346"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").Value = "'2.201"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'1.251"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "'1.186"
$ws.Range("E50").Value = "  +6.96%  "
$ws.Range("D51").Value = "'81.45"
$ws.Range("E51").Value = "  +0.26%  "
